$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F3 and G3 timestamps
$ws.Range("F3").Value = "2025-11-12 12:19:07"
$ws.Range("G3").Value = "2025-11-12 12:20:14"

# Clear H3 entirely (cell should no longer exist in the sheet)
$ws.Range("H3").ClearContents()
